# Rename the existing sheet from "Sheet1" to "largepart"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "largepart"

# Add a new worksheet "smallpart" after the first sheet
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "smallpart"

# Header row
$ws2.Range("A1").Value = "type"
$ws2.Range("B1").Value = "name"
$ws2.Range("C1").Value = "imgList1"
$ws2.Range("D1").Value = "imgList2"
$ws2.Range("E1").Value = "imgList3"

# Data rows: pig cuts
$pigUrl = "http://jkmeat.co.kr/skin_shop/standard/s_img/basic/JK-0414_M.jpg"
$cowUrl = "http://image.auction.co.kr/itemimage/1f/02/25/1f02252136.jpg"

$data = @(
    @("pig", "도가니살", $pigUrl),
    @("pig", "부채살",   $pigUrl),
    @("pig", "설깃살",   $pigUrl),
    @("pig", "항정살",   $pigUrl),
    @("pig", "목심살",   $pigUrl),
    @("cow", "채끝살",   $cowUrl),
    @("cow", "아롱사태", $cowUrl),
    @("cow", "삼각살",   $cowUrl),
    @("cow", "제비추리", $cowUrl),
    @("cow", "업진살",   $cowUrl)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Set column C width on smallpart sheet to match the original best-fit style
# (target XML width="57.42578125"; offset by the engine's internal +5/7 padding
# so the saved width lands as close as possible to the target value)
$ws2.Columns.Item(3).ColumnWidth = 56.711495535714285

# Update view/selection state to match target:
# largepart: selection C28, no longer scrolled/pinned to the bottom rows
$ws1.Range("C28").Select() | Out-Null

# smallpart becomes the active (selected) sheet/tab, with selection C11
$ws2.Activate()
$ws2.Range("C11").Select() | Out-Null
